$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the "Currency" row label to lowercase "currency"
$wsInput.Range("A6").Value = "currency"

# Update the currency value, dropping the trailing space, and restyle the cell
$wsInput.Range("B6").Value = "US Dollar"
$wsInput.Range("B6").Style = "Normal"
$wsInput.Range("B6").Interior.Color = 5296274

# Select A6:B6 on the input sheet
$wsInput.Range("A6:B6").Select()

# Make ProductLoanInput the active sheet/tab (was ProductLoanOutput)
$wsInput.Activate()

$wb.Save()
